$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the extra (empty) row at row 2, shifting rows 3-6 up by one.
# This moves the "TEST" data row from row 3 to row 2, and drops the
# now-unused last row (previously row 6), matching "all import in extra
# row remove".
$ws.Rows.Item(2).Delete()

# Update the active selection to match the post-edit state (row 2 selected).
$ws.Range("A2:XFD2").Select()
